$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "96.216.23"
Set-TextValue "E2" "  +0.61%  "

Set-TextValue "D3" "3.568.20"
Set-TextValue "E3" "  -0.86%  "

Set-TextValue "E4" "  -0.03%  "

Set-TextValue "D5" "240.63"
Set-TextValue "E5" "  +0.87%  "

Set-TextValue "D6" "653.14"
Set-TextValue "E6" "  -0.68%  "

Set-TextValue "D7" "1.61"
Set-TextValue "E7" "  +8.39%  "

Set-TextValue "E8" "  -0.54%  "

Set-TextValue "E9" "  +0.02%  "

Set-TextValue "D10" "1.05"
Set-TextValue "E10" "  +4.77%  "

Set-TextValue "D11" "3.568.56"
Set-TextValue "E11" "  -0.86%  "

Set-TextValue "D12" "43.18"
Set-TextValue "E12" "  +0.24%  "

Set-TextValue "E13" "  +0.81%  "

Set-TextValue "E14" "  +1.00%  "

Set-TextValue "D15" "4.231.12"
Set-TextValue "E15" "  -0.93%  "

Set-TextValue "D16" "96.059.06"
Set-TextValue "E16" "  +0.57%  "

Set-TextValue "D17" "0.0000258"
Set-TextValue "E17" "  +0.99%  "

Set-TextValue "D18" "3.569.49"
Set-TextValue "E18" "  -0.84%  "

Set-TextValue "D19" "7.76"
Set-TextValue "E19" "  +0.08%  "

Set-TextValue "D20" "12.52"
Set-TextValue "E20" "  -1.05%  "

Set-TextValue "D22" "0.523"
Set-TextValue "E22" "  +6.03%  "

Set-TextValue "E23" "  -5.84%  "

Set-TextValue "D24" "506.38"
Set-TextValue "E24" "  -0.80%  "

Set-TextValue "D25" "0.0000198"
Set-TextValue "E25" "  +1.50%  "

Set-TextValue "E26" "  +3.23%  "

Set-TextValue "D27" "95.88"
Set-TextValue "E27" "  -1.05%  "

Set-TextValue "D28" "12.59"
Set-TextValue "E28" "  -1.16%  "

Set-TextValue "D29" "3.760.40"
Set-TextValue "E29" "  -0.90%  "

Set-TextValue "E30" "  +7.45%  "

Set-TextValue "D31" "2.97"
Set-TextValue "E31" "  -6.43%  "

Set-TextValue "D32" "11.36"
Set-TextValue "E32" "  +0.29%  "

Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  -0.01%  "

Set-TextValue "E34" "  +3.05%  "

Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.31%  "

Set-TextValue "D36" "31.36"

Set-TextValue "D37" "616.97"
Set-TextValue "E37" "  +7.20%  "

Set-TextValue "E38" "  +6.99%  "

Set-TextValue "D39" "0.562"
Set-TextValue "E39" "  +0.32%  "

Set-TextValue "D40" "1.61"
Set-TextValue "E40" "  +8.61%  "

Set-TextValue "E41" "  +0.05%  "

Set-TextValue "E42" "  -0.30%  "

Set-TextValue "D43" "0.900"
Set-TextValue "E43" "  -2.53%  "

Set-TextValue "E44" "  +5.58%  "

Set-TextValue "D50" "3.53"
Set-TextValue "E50" "  +1.43%  "

Set-TextValue "E51" "  +0.32%  "

# Reordered coin rows 45-49
Set-TextValue "B45" "Filecoin"
Set-TextValue "C45" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D45" "5.68"
Set-TextValue "E45" "  -1.08%  "

Set-TextValue "B46" "WhiteBITCoin"
Set-TextValue "C46" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D46" "23.50"
Set-TextValue "E46" "  -1.16%  "

Set-TextValue "B47" "Stacks"
Set-TextValue "C47" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "2.27"
Set-TextValue "E47" "  +1.23%  "

Set-TextValue "B48" "VeChain"
Set-TextValue "C48" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0419"
Set-TextValue "E48" "  +0.73%  "

Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "33.42"
Set-TextValue "E49" "  -1.16%  "
